$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample data")

# Update sx/sy values in the table
$ws.Range("B5").Value = 5.07
$ws.Range("C5").Value = 8.3
$ws.Range("B6").Value = 5.07
$ws.Range("C6").Value = 8.3
$ws.Range("C9").Value = 6

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("C6").Select()
